$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 8: APMS-T143 / Email = awtindia.sc@gmail.com ----------------------
$ws.Range("A8").Value = "APMS-T143"

# --- New column D header ---------------------------------------------------
$ws.Range("D1").Value = "Email"

$ws.Range("D8").Value = "awtindia.sc@gmail.com"

# --- Row 9: APMS-T144 / Email = nicoalastestla458@gmail.com (hyperlinked) --
$ws.Range("A9").Value = "APMS-T144"
$ws.Range("D9").Value = "nicoalastestla458@gmail.com"
$ws.Hyperlinks.Add($ws.Range("D9"), "mailto:nicoalastestla458@gmail.com")

# --- Match formatting of the existing table ---------------------------------
# Column A (TestCaseId) cells use the centered "Normal" style from A2.
$ws.Range("A2").Copy()
$ws.Range("A8:A9").PasteSpecial(-4122)

# D9 is a hyperlink cell like C6/C7 - copy that style (applied after the
# hyperlink was added so the Hyperlink style sticks).
$ws.Range("C6").Copy()
$ws.Range("D9").PasteSpecial(-4122)

$ws.Application.CutCopyMode = $false

# New column D width
$ws.Range("D1").ColumnWidth = 28

# --- View state --------------------------------------------------------
$ws.Range("D11").Select() | Out-Null
